$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Duty Cycle" - elaborate duty-cycle calculus filled in (Vd, Rdson,
# Iout worst-case values) and the "Duty Cycle approximation difference"
# result moved from H11 into the merged E11:G11 block (G11).
# ---------------------------------------------------------------------------
$wsDuty = $wb.Worksheets.Item("Duty Cycle")

$wsDuty.Range("F5").Value = 0.4
$wsDuty.Range("F6").Value = 0.195
$wsDuty.Range("G6").Value = "Worst case"
$wsDuty.Range("F7").Value = 0.6
$wsDuty.Range("G7").Value = "Worst case"
$wsDuty.Range("G8").Value = "Voltage dropout accross the internal MOSFET"

$diffFormula = $wsDuty.Range("H11").Formula
$wsDuty.Range("E11:G11").UnMerge()
$wsDuty.Range("G11").Formula = $diffFormula
$wsDuty.Range("G11").NumberFormat = "0%"
$wsDuty.Range("H11").Clear()
$wsDuty.Range("E11:F11").Merge()

# ---------------------------------------------------------------------------
# Sheet "Inductor value" - ripple target raised from 10% to 17%, and the
# inductance formula now uses the elaborate duty cycle (F9) instead of the
# simple one (C5).
# ---------------------------------------------------------------------------
$wsInd = $wb.Worksheets.Item("Inductor value")

$wsInd.Range("C9").Value = 0.17
$wsInd.Range("C12").Formula = "=('Duty Cycle'!F9*C11)/(2*D9)*('Duty Cycle'!C3-'Duty Cycle'!C4)"

# ---------------------------------------------------------------------------
# Sheet "Capacitors" - Irms formula now uses the elaborate duty cycle (F9),
# and a new "in the other way around" worst-case ΔVout block is added.
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacitors")

$wsCap.Range("C11").Formula = "='Inductor value'!C8*SQRT('Duty Cycle'!F9*(1-'Duty Cycle'!F9))"

$wsCap.Range("I19").Value = "Resr (Ω)"
$wsCap.Range("I20").Value = "Cout (F)"
$wsCap.Range("J20").Value = 0.000022
$wsCap.Range("I21").Value = "ΔVout (V)"
$wsCap.Range("J21").Formula = "='Inductor value'!D9*(J19+1/(8*'Inductor value'!C10*J20))"
$wsCap.Range("L19").Value = "In the other way around"

# ---------------------------------------------------------------------------
# Sheet "Diode" - Id1 formula now uses the elaborate duty cycle (F9).
# ---------------------------------------------------------------------------
$wsDiode = $wb.Worksheets.Item("Diode")

$wsDiode.Range("C13").Formula = "=('Inductor value'!C8)*(1-'Duty Cycle'!F9)"

# ---------------------------------------------------------------------------
# Sheet "Feedback resistors" - no formula/value changes.
# ---------------------------------------------------------------------------
$wsFb = $wb.Worksheets.Item("Feedback resistors")

# ---------------------------------------------------------------------------
# Selections / active sheet, matching the final view state. Selecting a
# range implicitly activates its sheet, so the order here matters: the
# last sheet selected/activated ends up the active tab ("Feedback
# resistors", matching workbook activeTab=4).
# ---------------------------------------------------------------------------
$wsDuty.Range("E11:G11").Select()
$wsInd.Range("C10").Select()
$wsCap.Range("L18").Select()
$wsDiode.Range("C14").Select()
$wsFb.Range("C16").Select()
